# Apply changes to pkm_pokedexes workbook:
# 1. Column F (is_main_series): convert numeric 0/1 values to text "TRUE"/"FALSE"
# 2. Fix apostrophes in three description cells (E7, E8, E15): ' -> ''

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix doubled apostrophes in description text (column E) first, so the
#     shared-string table gets these new strings added before TRUE/FALSE ---
$ws.Range("E7").Value = "'Platinum Sinnoh dex—an extended version of Diamond and Pearl''s"
$ws.Range("E8").Value = "'HeartGold/SoulSilver Johto dex—Gold/Silver/Crystal''s, extended to add move-based Generation IV evolutions"
$ws.Range("E15").Value = "'Omega Ruby/Alpha Sapphire Hoenn Dex—Ruby/Sapphire/Emerald''s, updated to add new evolutions"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("E15").Style = "Normal"

# --- Update column F (is_main_series) for rows 2-15: numeric 0/1 -> text "TRUE"/"FALSE" ---
$trueRows = @(2,3,4,5,6,7,8,9,10,12,13,14,15)
$falseRows = @(11)

foreach ($r in $trueRows) {
    $ws.Cells.Item($r, 6).Value = "'TRUE"
}
foreach ($r in $falseRows) {
    $ws.Cells.Item($r, 6).Value = "'FALSE"
}

# Leading apostrophe forces Excel to store the value as literal text instead
# of auto-converting to a Boolean; reset the cell style afterwards so the
# quote-prefix indicator isn't left applied to the cells.
$ws.Range("F2:F15").Style = "Normal"

Write-Host "Done applying changes"
